# Annot 30 -> 50
# Fill in the "label" (annotation) column E for rows 32 through 52 of Sheet1.
# Rows with "?" mean "no clear winner / undecided"; numeric rows are the
# annotator's preference score.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    32 = "?"
    33 = "?"
    34 = -2
    35 = -1
    36 = -1
    37 = -2
    38 = 0
    39 = 0
    40 = "?"
    41 = 0
    42 = -2
    43 = 0
    44 = "?"
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = -2
    50 = -2
    51 = 0
    52 = 0
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Range("E$row").Value = $values[$row]
}

# Move the view/selection to reflect where the user ended up editing.
$ws.Activate()
$ws.Range("E53").Select()
